# Apply the "updated trading sheet and other models" edit:
#  - Append 10 new telecom/media rows (34-43) to the Main sheet
#  - Update window/view/zoom state saved with the workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append: Name, Ticker pairs
$newRows = @(
    @("Singapore Telecom", "ST SP"),
    @("Telstra", "TLS AU"),
    @("Telus", "T CN"),
    @("Cellnex", "CLNX"),
    @("Orange", "ORA FP"),
    @("Swisscom", "SCMN SW"),
    @("Telekom Indonesia", "TLKM IJ"),
    @("Telefonica", "TEF SM"),
    @("Wolters Kluwer", "WKL NA"),
    @("Sirius XM", "SIRI")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $name = $newRows[$i][0]
    $ticker = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = "x"
    $ws.Cells.Item($r, 2).Value = $name
    $ws.Cells.Item($r, 3).Value = $ticker
}

# Update the sheet view: zoom to 175%, scroll frozen pane down, move selection
$ws.Application.ActiveWindow.Zoom = 175
$ws.Activate()
$ws.Range("C36").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A44").Select()

# Update workbook window placement/size to match the author's last save
$excel.Width = 24495
$excel.Height = 16830
$excel.Left = 53535
$excel.Top = 855
